$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item(3)

# --- Header row styling (A1/B1 become bold, matching style index 6) ---
$ws3.Range("A1:B1").Font.Bold = $true

# --- New rows 33-37 (Rationing Program block) ---
$ws3.Range("A33").Value = "Storage Volume Trigger for Rationing Programs (AF)"
$ws3.Range("B33").Value = 1000000
$ws3.Range("B33").NumberFormat = "#,##0.00"

$ws3.Range("A34").Value = "Cost for Rationing Program (`$/capita)"
$ws3.Range("B34").Value = 2

$ws3.Range("A35").Value = "Consecutive Year Loss Adjustment (%)"
$ws3.Range("B35").Value = 0

$ws3.Range("A36").Value = "Demand Hardening Factor (%)"
$ws3.Range("B36").Value = 50

$ws3.Range("A37").Value = "Retail Price (`$/AF)"
$ws3.Range("B37").Value = 650

# --- Header for rationing program block (row 32) ---
$ws3.Range("A32").Value = "contingentWMOsInput_rationingProgram.csv"
$ws3.Range("A32").Font.Bold = $true

# --- Header for elasticity of demand block (row 39) ---
$ws3.Range("A39").Value = "contingentWMOsInput_elasticityofDemand.csv"
$ws3.Range("A39").Font.Bold = $true

# --- Elasticity of demand rows ---
$ws3.Range("A40").Value = "Elasticity of Demand Single Family Residential"
$ws3.Range("B40").Value = -0.2

$ws3.Range("A42").Value = "Elasticity of Demand Industrial"
$ws3.Range("B42").Value = -0.1

$ws3.Range("A43").Value = "Elasticity of Demand Commercial and Governmental"
$ws3.Range("B43").Value = -0.11

$ws3.Range("A44").Value = "Elasticity of Demand Landscape"
$ws3.Range("B44").Value = -0.4

$ws3.Range("A45").Value = "Lower Loss Boundary"
$ws3.Range("B45").Value = 0

$ws3.Range("A46").Value = "Upper Loss Boundary"
$ws3.Range("B46").Value = 0.7

$ws3.Range("A41").Value = "Elasticity of Demand Multi-Family Residential"
$ws3.Range("B41").Value = -0.12

# --- Column B width on TestInputData ---
$ws3.Columns.Item(2).ColumnWidth = 21.166666666666668

# --- Sheet2 (test_modelLogic.py): move selection, remove tab-selected state ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Activate()
$ws2.Range("A9").Select()

# --- Sheet3 (TestInputData) becomes the active/selected tab ---
$ws3.Activate()
$ws3.Range("A42").Select()

Write-Host "done"
